$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(30, 6.24, 109.193),
    @(31, 5.26, 115.1801),
    @(32, 6.26, 120.8746),
    @(33, 6.27, 126.6955),
    @(34, 6.28, 132.71549999999999),
    @(35, 7.28, 139.41759999999999)
)

$ws.Range("C30").Copy() | Out-Null
$ws.Range("C31:C36").PasteSpecial(-4122) | Out-Null

$row = 31
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}

$ws.Range("D16").Select()
